# Update the 丽水-漫展信息 workbook: refresh the event listing on the
# "展览" and "全部类型" sheets. Two outdated events (2024-08-17 entries)
# drop off the list, the two upcoming events move up into rows 2-3 with
# refreshed data, and the now-unused rows 4-5 are removed.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    # Row 2 becomes the "丽水·Re动漫游戏嘉年华" event (previously row 4),
    # with the attendee count refreshed from 507 to 509.
    # Column B holds "yyyy-mm-dd" text; force text formatting first so the
    # COM layer doesn't silently reinterpret it as a date serial, then drop
    # the temporary number format again so the cell keeps its original
    # (default) style.
    $ws.Cells.Item(2, 2).NumberFormat = "@"
    $ws.Cells.Item(2, 2).Value = "2024-08-24"
    $ws.Cells.Item(2, 2).ClearFormats()
    $ws.Cells.Item(2, 3).Value = "丽水·Re动漫游戏嘉年华"
    $ws.Cells.Item(2, 4).Value = "中东路848号(解放街交汇) 飞达国际大酒店"
    $ws.Cells.Item(2, 5).Value = "2024.08.24 09:30-08.24 17:00"
    $ws.Cells.Item(2, 6).Value = 509
    $ws.Cells.Item(2, 7).Value = 45
    $ws.Cells.Item(2, 8).Value = "https://show.bilibili.com/platform/detail.html?id=89651"
    $ws.Cells.Item(2, 9).Value = "//i0.hdslb.com/bfs/openplatform/202407/7o5ALbAM1721383424201.jpeg"

    # Row 3 becomes the "丽水·LZ栗子动漫游戏嘉年华" event (previously row 5).
    $ws.Cells.Item(3, 2).NumberFormat = "@"
    $ws.Cells.Item(3, 2).Value = "2024-09-16"
    $ws.Cells.Item(3, 2).ClearFormats()
    $ws.Cells.Item(3, 3).Value = "丽水·LZ栗子动漫游戏嘉年华"
    $ws.Cells.Item(3, 4).Value = "城北街798号 莱茵体育生活馆"
    $ws.Cells.Item(3, 5).Value = "2024.09.16 09:30-09.16 17:00"
    $ws.Cells.Item(3, 6).Value = 444
    $ws.Cells.Item(3, 7).Value = 65
    $ws.Cells.Item(3, 8).Value = "https://show.bilibili.com/platform/detail.html?id=87480"
    $ws.Cells.Item(3, 9).Value = "//i1.hdslb.com/bfs/openplatform/202406/bATqcZhH1719285865931.jpeg"

    # The old rows 4 and 5 (now duplicated into rows 2-3 above) are removed,
    # shrinking the sheet from A1:I5 to A1:I3.
    $ws.Rows("4:5").Delete()
}
